$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new trade row (row 8) below the existing data, following the
# same layout/formatting used for the preceding rows (3:7).

# Column A holds a date/time serial; copy row 7's number format down to
# row 8 first so the new cell renders the same way as the rows above it.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A8").Value = 42650.371516203704   # Date
$ws.Range("B8").Value = $false               # Profitable
$ws.Range("C8").Value = 9909.2000000000007   # Principle
$ws.Range("D8").Value = 9994.15              # Start Principle
$ws.Range("E8").Value = 104.839996           # BuyPrice
$ws.Range("F8").Value = 103.95               # SellPrice
$ws.Range("G8").Value = $false               # IsShortSell
$ws.Range("H8").Value = -0.85                # Price Change %
$ws.Range("I8").Value = $true                # Strong trade

# The new BuyPrice value (104.839996) is wider than anything previously in
# column E, so its best-fit width grows to accommodate it.
$ws.Columns.Item(5).ColumnWidth = 10
